$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (New York -- New York): clear B4:H4, K4:L4, set J4 = FALSE, update O4
$ws.Range("B4:H4").Clear()
$ws.Range("K4:L4").Clear()
$ws.Range("J4").Value = $false
$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"

# Row 39 (Delaware): update O39
$ws.Range("O39").Value = "An error occurred. ... AttributeError(""'numpy.float64' object has no attribute 'split'"")"
